# Refresh the "cryptos" price/volume table with the latest scrape.
#
# D/E are stored as plain text in the workbook (prices use "."
# as a thousands separator, e.g. "68.544.97", and percentages keep
# padding spaces, e.g. "  +3.06%  "), never as real numbers/percents.
# For values that *look* like an ordinary decimal (e.g. "10.15"),
# assigning .Value directly would make Excel silently reinterpret
# them as numbers (and drop meaningful trailing zeros, e.g. "12.60"
# -> 12.6). To keep them as text we prefix with a literal leading
# apostrophe (Excel's own 'treat as text' convention) and then
# clear the resulting quote-prefix formatting so no stray cell style
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.544.97'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").Value = '3.653.18'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '''199.12'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +9.79%  '
$ws.Range("D6").Value = '''578.35'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("D7").Value = '3.649.36'
$ws.Range("E7").Value = '  +2.50%  '
$ws.Range("E8").Value = '  +2.13%  '
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("E11").Value = '  +8.34%  '
$ws.Range("D12").Value = '''56.73'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +5.80%  '
$ws.Range("D13").Value = '''0.0000297'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +18.10%  '
$ws.Range("D14").Value = '''10.15'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.25%  '
$ws.Range("D15").Value = '4.237.16'
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").Value = '3.650.91'
$ws.Range("E16").Value = '  +2.62%  '
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").Value = '''12.61'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.46%  '
$ws.Range("D19").Value = '68.493.20'
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("D20").Value = '''18.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("E21").Value = '  +4.23%  '
$ws.Range("D22").Value = '''405.57'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.50%  '
$ws.Range("D23").Value = '''13.18'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +29.12%  '
$ws.Range("E24").Value = '  -0.32%  '
$ws.Range("D25").Value = '''86.15'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.18%  '
$ws.Range("E26").Value = '  +3.79%  '
$ws.Range("D27").Value = '''12.71'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.71%  '
$ws.Range("E28").Value = '  +8.88%  '
$ws.Range("D29").Value = '''6.13'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.65%  '
$ws.Range("D30").Value = '''8.29'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +23.53%  '
$ws.Range("D31").Value = '''9.27'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.44%  '
$ws.Range("D32").Value = '''32.12'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.74%  '
$ws.Range("D33").Value = '''696.15'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +15.40%  '
$ws.Range("E34").Value = '  +3.59%  '
$ws.Range("D35").Value = '''0.117'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.45%  '
$ws.Range("D36").Value = '''64.81'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("D37").Value = '''42.90'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.95%  '
$ws.Range("D38").Value = '''0.430'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +15.90%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  +8.04%  '
$ws.Range("E41").Value = '  +10.81%  '
$ws.Range("D42").Value = '''2.92'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +21.45%  '
$ws.Range("D43").Value = '''3.17'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +14.66%  '
$ws.Range("D44").Value = '3.212.08'
$ws.Range("E44").Value = '  +10.72%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''0.999'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '''3.01'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +36.29%  '
$ws.Range("E47").Value = '  +4.59%  '
$ws.Range("E48").Value = '  +9.38%  '
$ws.Range("E49").Value = '  +2.38%  '
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("E51").Value = '  +5.87%  '
